# Auto-generated edit script applying scheduled-runner price updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 21116.877
$ws.Range("I132").Value = 3229.0527
$ws.Range("J132").Value = 82911.17999999999
$ws.Range("K132").Value = 9687.158100000001
$ws.Range("L132").Value = 248733.54
$ws.Range("M132").Value = -7157.158100000001
$ws.Range("N132").Value = -253793.54

$ws.Range("H133").Value = 56263.168
$ws.Range("J133").Value = 56263.168
$ws.Range("L133").Value = 56263.168
$ws.Range("N133").Value = -66383.16800000001

$ws.Range("H138").Value = 1446.06
$ws.Range("J138").Value = 1964.6833
$ws.Range("L138").Value = 5894.0499
$ws.Range("N138").Value = -16174.0499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20340

$ws.Range("H19").Value = 13000
$ws.Range("J19").Value = 13000
$ws.Range("L19").Value = 13000
$ws.Range("N19").Value = -13458

$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20540

$ws.Range("H45").Value = 2494.4443
$ws.Range("J45").Value = 2588
$ws.Range("L45").Value = 2588
$ws.Range("N45").Value = -3342

$ws.Range("H74").Value = 1524.6123
$ws.Range("I74").Value = 1291.091
$ws.Range("K74").Value = 1291.091
$ws.Range("M74").Value = -417.0909999999999

$ws.Range("H77").Value = 1524.6123
$ws.Range("I77").Value = 1291.091
$ws.Range("K77").Value = 6455.455
$ws.Range("M77").Value = -2087.455

$ws.Range("H98").Value = 39400
$ws.Range("J98").Value = 39400
$ws.Range("L98").Value = 39400
$ws.Range("N98").Value = -45390

$ws.Range("H102").Value = 52654.75
$ws.Range("I102").Value = 3533
$ws.Range("J102").Value = 200020
$ws.Range("K102").Value = 3533
$ws.Range("L102").Value = 200020
$ws.Range("M102").Value = -1911
$ws.Range("N102").Value = -203264

$ws.Range("H109").Value = 44876.5
$ws.Range("J109").Value = 44876.5
$ws.Range("L109").Value = 44876.5
$ws.Range("N109").Value = -47650.5

$ws.Range("H132").Value = 2662.56
$ws.Range("I132").Value = 1492.6923
$ws.Range("J132").Value = 3929.9167
$ws.Range("K132").Value = 4478.0769
$ws.Range("L132").Value = 11789.7501
$ws.Range("M132").Value = -1948.0769
$ws.Range("N132").Value = -16849.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H134").Value = 3532.7407
$ws.Range("I134").Value = 4288.75
$ws.Range("K134").Value = 12866.25
$ws.Range("M134").Value = -10331.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 39500
$ws.Range("J28").Value = 39500
$ws.Range("L28").Value = 39500
$ws.Range("N28").Value = -39990

$ws.Range("H95").Value = 90000
$ws.Range("J95").Value = 90000
$ws.Range("L95").Value = 90000
$ws.Range("N95").Value = -95492

$ws.Range("H122").Value = 100899.914
$ws.Range("I122").Value = 240799.8
$ws.Range("J122").Value = 971.4286
$ws.Range("K122").Value = 722399.3999999999
$ws.Range("L122").Value = 2914.2858
$ws.Range("M122").Value = -719949.3999999999
$ws.Range("N122").Value = -7814.2858

$ws.Range("H132").Value = 112516.08
$ws.Range("I132").Value = 2820
$ws.Range("J132").Value = 181076.12
$ws.Range("K132").Value = 8460
$ws.Range("L132").Value = 543228.36
$ws.Range("M132").Value = -5930
$ws.Range("N132").Value = -548288.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5529.32
$ws.Range("I5").Value = 6041.5557
$ws.Range("J5").Value = 4212.143
$ws.Range("K5").Value = 18124.6671
$ws.Range("L5").Value = 12636.429
$ws.Range("M5").Value = -18012.6671
$ws.Range("N5").Value = -12860.429

$ws.Range("H9").Value = 71435460
$ws.Range("I9").Value = 500005000
$ws.Range("K9").Value = 1500015000
$ws.Range("M9").Value = -1500014776

$ws.Range("H12").Value = 320.64517
$ws.Range("I12").Value = 245.42857
$ws.Range("J12").Value = 342.58334
$ws.Range("K12").Value = 736.28571
$ws.Range("L12").Value = 1027.75002
$ws.Range("M12").Value = -563.28571
$ws.Range("N12").Value = -1373.75002

$ws.Range("H135").Value = 5529.32
$ws.Range("I135").Value = 6041.5557
$ws.Range("J135").Value = 4212.143
$ws.Range("K135").Value = 54374.0013
$ws.Range("L135").Value = 37909.287
$ws.Range("M135").Value = -51839.0013
$ws.Range("N135").Value = -42979.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 21975
$ws.Range("J62").Value = 22000
$ws.Range("L62").Value = 22000
$ws.Range("N62").Value = -23372

$ws.Range("H65").Value = 21975
$ws.Range("J65").Value = 22000
$ws.Range("L65").Value = 66000
$ws.Range("N65").Value = -72864

$ws.Range("H100").Value = 38355
$ws.Range("J100").Value = 38355
$ws.Range("L100").Value = 38355
$ws.Range("N100").Value = -40519

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3291.2727
$ws.Range("I7").Value = 2680.8
$ws.Range("J7").Value = 3800
$ws.Range("K7").Value = 2680.8
$ws.Range("L7").Value = 3800
$ws.Range("M7").Value = -2568.8
$ws.Range("N7").Value = -4024

$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21498

$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -67488

$ws.Range("H126").Value = 3291.2727
$ws.Range("I126").Value = 2680.8
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 8042.400000000001
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -5572.400000000001
$ws.Range("N126").Value = -16340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1433.76
$ws.Range("I132").Value = 1116.1945
$ws.Range("J132").Value = 2250.3572
$ws.Range("K132").Value = 3348.5835
$ws.Range("L132").Value = 6751.071599999999
$ws.Range("M132").Value = -818.5835000000002
$ws.Range("N132").Value = -11811.0716
